# Generate Report for Handoff
# Update status "In Translation" -> "Ready for handoff" and refresh the
# handoff timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 17:13:07"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2).
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 17:12:57"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 17:13:07"

# The longer "Ready for handoff" text no longer fits the previous column
# width, so re-autofit the affected columns (matches Excel's automatic
# column-width recalculation on edit).
$overview.Range("E:E").ColumnWidth = 16.333333333333332
$overview.Range("F:F").ColumnWidth = 16.333333333333332
$zhcn.Range("C:C").ColumnWidth = 16.333333333333332
$dede.Range("C:C").ColumnWidth = 16.333333333333332
